$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 24 (shifts existing rows 24..80 down to 25..81),
# inheriting number formats / styles from the row above (row 23), which already
# matches the style pattern we need (s="1" for A-F, s="2" for G, s="3" for I/J).
$ws.Rows("24:24").Insert()

# --- Fill in the new row's data: "0199 - LEI" (Legal Entity Identifier) ---
# Values are assigned in the same order the original author entered them
# (B, D, E, J, I, A) so newly-created shared-string entries land at the same
# indices as the canonical file; C/F reuse pre-existing shared strings either way.
$ws.Range("B24").Value = "0199"
$ws.Range("D24").Value = "Legal Entity Identifier (LEI)"
$ws.Range("E24").Value = "As of December 2018, there are 33 LEI issuing organizations in the world."
$ws.Range("J24").Value = "The entire 20 character code (including the check digits)"
$ws.Range("I24").Value = "The ISO 17442 standard specifies the minimum reference data, which must the format of the organization identifiers, be Supplied for each LEI:`r`n* The official name of the legal entity as recorded in the official registers.`r`n* The registered address of that legal entity.`r`n* The country of formation.`r`n* The codes for the representation of names of countries and their subdivisions.`r`n* The date of the first LEI assignment; the date of last update of the`r`n* LEI information; and the date of expiry, if applicable.`r`nAdditional information may be registered as agreed between the legal entity and its LEI issuing organization.`r`netc. - see ICD sheet for further information"
$ws.Range("A24").Value = "LEI"
$ws.Range("C24").Value = "international"
$ws.Range("F24").Value = "5"
$ws.Range("G24").Formula = "=FALSE"

# Row height for the new row (matches the target layout after the text was added)
$ws.Rows("24:24").RowHeight = 174

# Refresh the AutoFilter so its range grows to match the new used range (A1:L81).
# Toggling off first avoids AutoFilter() simply disabling the existing filter.
$ws.AutoFilterMode = $false
$ws.Range("A1:L81").AutoFilter()

# Keep the hidden _FilterDatabase defined name's range in step with the filter.
$wb.Names.Item(1).RefersTo = "='Participant Identifier Scheme'!`$A`$1:`$L`$81"

# Restore the selection to the new row, matching where the editor ended up.
$ws.Range("A24:D24").Select()
